$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3C")

# C7 was stored as text "31"; normalize it to a real number.
$ws.Range("C7").Value = 31

# Append the new submission row synced at 2026-02-10 14:21:14.
$ws.Range("A8").Value = "2026-02-10 14:21:14"
$ws.Range("B8").Value = "FATIMA ALHAJI GANA"
$ws.Range("C8").Value = 45
$ws.Range("D8").Value = 9
